$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.734.07'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '1.754.83'
$ws.Range("E3").Value = '  -2.84%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'" + '324.57'
$ws.Range("E5").Value = '  -4.10%  '
$ws.Range("D6").Value = "'" + '0.9989'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = "'" + '0.4288'
$ws.Range("E7").Value = '  -7.81%  '
$ws.Range("D8").Value = "'" + '0.3649'
$ws.Range("E8").Value = '  -4.21%  '
$ws.Range("D9").Value = "'" + '45.38'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").Value = "'" + '0.07495'
$ws.Range("E10").Value = '  -1.29%  '
$ws.Range("E11").Value = '  -3.15%  '
$ws.Range("D12").Value = "'" + '0.9996'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = "'" + '21.69'
$ws.Range("E13").Value = '  -3.54%  '
$ws.Range("D14").Value = "'" + '6.156'
$ws.Range("E14").Value = '  -3.19%  '
$ws.Range("D15").Value = "'" + '7.266'
$ws.Range("E15").Value = '  -3.41%  '
$ws.Range("D16").Value = '1.747.71'
$ws.Range("E16").Value = '  -3.38%  '
$ws.Range("D17").Value = "'" + '0.00001070'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("D18").Value = "'" + '87.89'
$ws.Range("E18").Value = '  +7.68%  '
$ws.Range("D19").Value = "'" + '0.06207'
$ws.Range("E19").Value = '  -7.80%  '
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = "'" + '17.12'
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("D22").Value = "'" + '6.163'
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").Value = "'" + '0.5269'
$ws.Range("E23").Value = '  -4.87%  '
$ws.Range("D24").Value = '27.736.38'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("D25").Value = "'" + '11.71'
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("D26").Value = "'" + '2.335'
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("D27").Value = "'" + '20.58'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Value = "'" + '152.91'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("D29").Value = "'" + '2.370'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Value = '1.949.17'
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("D31").Value = "'" + '1.226'
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").Value = "'" + '127.60'
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("D33").Value = "'" + '5.751'
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").Value = "'" + '0.09160'
$ws.Range("E34").Value = '  -4.91%  '
$ws.Range("D35").Value = "'" + '3.654'
$ws.Range("E35").Value = '  -9.54%  '
$ws.Range("D36").Value = "'" + '12.76'
$ws.Range("E36").Value = '  +5.39%  '
$ws.Range("D37").Value = "'" + '0.02318'
$ws.Range("D38").Value = "'" + '0.2156'
$ws.Range("E38").Value = '  -7.69%  '
$ws.Range("D39").Value = "'" + '5.130'
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("D41").Value = "'" + '0.06116'
$ws.Range("E41").Value = '  -3.93%  '
$ws.Range("D42").Value = "'" + '1.198'
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("D43").Value = "'" + '1.428'
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("D44").Value = "'" + '7.999'
$ws.Range("E44").Value = '  -4.43%  '
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").Value = "'" + '13.76'
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("D47").Value = "'" + '0.5961'
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("D48").Value = "'" + '3.753'
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("D49").Value = "'" + '126.36'
$ws.Range("E49").Value = '  -3.77%  '
$ws.Range("D50").Value = "'" + '1.976'
$ws.Range("E50").Value = '  -3.55%  '
$ws.Range("D51").Value = "'" + '0.06907'
$ws.Range("E51").Value = '  -3.47%  '

# Reset style on text-forced cells so the quote-prefix marker
# introduced above does not leave a stray cell style behind.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
